# Actualización automática 2025-09-01 08:30:07
$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" ---
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

$wsGrupo.Range("H13").Value = 0
$wsGrupo.Range("I13").Value = 0
$wsGrupo.Range("L18").Value = 0
$wsGrupo.Range("L19").Value = 0
$wsGrupo.Range("M19").Value = 0

$wsGrupo.Range("H29").Value = "0 de 27"
$wsGrupo.Range("I29").Value = "0 de 27"
$wsGrupo.Range("L29").Value = "0 de 27"
$wsGrupo.Range("M29").Value = "0 de 27"

# --- Sheet "VENTA MENSUAL" ---
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

# Rolling month headers: mayo/junio/julio/agosto -> junio/julio/agosto/septiembre
$wsMensual.Range("C1").Value = "junio"
$wsMensual.Range("D1").Value = "julio"
$wsMensual.Range("E1").Value = "agosto"
$wsMensual.Range("F1").Value = "septiembre"

# Column F widened to fit "septiembre" (stored width 16 == ColumnWidth 16 - 5/6)
$wsMensual.Columns.Item(6).ColumnWidth = 15.166666666666666

# Row 4 - ARBOLEDA ZAMBRANO ROBERTO ANTONIO
$wsMensual.Range("C4").Value = 0

# Row 13 - GRANIMUNDO S.A.
$wsMensual.Range("C13").Value = 738.66
$wsMensual.Range("D13").Value = 0
$wsMensual.Range("E13").Value = 1284.64
$wsMensual.Range("F13").Value = 0

# Row 14 - LINO TUMBACO VICENTE JAVIER
$wsMensual.Range("C14").Value = 1473.73
$wsMensual.Range("D14").Value = 3990.41
$wsMensual.Range("E14").Value = 0

# Row 17 - MOREIRA MOREIRA PATRICIO IGNACIO
$wsMensual.Range("C17").Value = 9556.26
$wsMensual.Range("D17").Value = 0

# Row 18 - PAREDES ORTIZ MARIA INES
$wsMensual.Range("C18").Value = 64.82
$wsMensual.Range("D18").Value = 8691.84
$wsMensual.Range("E18").Value = 4413.66
$wsMensual.Range("F18").Value = 0

# Row 19 - RENOVA&DISEÑA S.A.
$wsMensual.Range("C19").Value = 411.7
$wsMensual.Range("D19").Value = 2045.31
$wsMensual.Range("E19").Value = 1837.88
$wsMensual.Range("F19").Value = 0

# Row 21 - ROCA REYNA PAUL DAVID
$wsMensual.Range("C21").Value = 3225.33
$wsMensual.Range("D21").Value = 0

# Row 29 - totals
$wsMensual.Range("C29").Value = 15470.5
$wsMensual.Range("D29").Value = 14727.56
$wsMensual.Range("E29").Value = 7536.18
$wsMensual.Range("F29").Value = 0
